$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the "theta_threshold_range" row (row 5), shifting row 6 up to row 5
$ws.Rows.Item(5).Delete()

# Update values for rows 2-5 (B and C columns)
$ws.Range("B2").Value = 5.4
$ws.Range("C2").Value = 10
$ws.Range("B3").Value = 5
$ws.Range("C3").Value = 8.9
$ws.Range("B4").Value = 0.9
$ws.Range("C4").Value = 1.2
$ws.Range("B5").Value = 0
$ws.Range("C5").Value = 15

# Fix style of C5 (was using a different font style before the shift)
$ws.Range("C5").Style = $ws.Range("B5").Style

# Set selection
$ws.Range("C5").Select()

$wb.Windows.Item(1).WindowState = -4143
